# ninjacrm testData.xlsx update:
#   - rename the sheet from "Sheet1" to "Opp"
#   - re-key row 2 (OPP-1000) so Amount*, Expected Close Date and
#     Probability are stored as quote-prefixed text instead of numbers,
#     matching how the rest of the sheet was re-entered by hand
#   - move the active selection to F2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Opp"

# Amount* (C2) and Probability (I2): force text storage via a leading
# apostrophe, same as Excel does when a user types '26847 / '64.
$ws.Range("C2").Value = "'26847"
$ws.Range("I2").Value = "'64"

# Expected Close Date (F2): replace the ISO date with a quote-prefixed
# text date (04-06-2025), formatted with the short-date display.
$ws.Range("F2").Value = "'04-06-2025"
$ws.Range("F2").NumberFormat = "mm-dd-yy"

$ws.Range("F2").Select()
